$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4416279.5
$ws.Range("J17").Value = 4800277.5
$ws.Range("L17").Value = 14400832.5
$ws.Range("N17").Value = -14401168.5
$ws.Range("H32").Value = 1343.375
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2674
$ws.Range("H93").Value = 35601
$ws.Range("J93").Value = 35601
$ws.Range("L93").Value = 35601
$ws.Range("N93").Value = -40593
$ws.Range("H98").Value = 2405.5625
$ws.Range("I98").Value = 2366.0667
$ws.Range("K98").Value = 2366.0667
$ws.Range("M98").Value = -868.0666999999999
$ws.Range("H122").Value = 2405.5625
$ws.Range("I122").Value = 2366.0667
$ws.Range("K122").Value = 7098.2001
$ws.Range("M122").Value = -4648.2001
$ws.Range("H132").Value = 3187.282
$ws.Range("I132").Value = 3018.625
$ws.Range("J132").Value = 3958.2856
$ws.Range("K132").Value = 9055.875
$ws.Range("L132").Value = 11874.8568
$ws.Range("M132").Value = -6525.875
$ws.Range("N132").Value = -16934.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9140.769
$ws.Range("I32").Value = 8193.549999999999
$ws.Range("J32").Value = 15455.556
$ws.Range("K32").Value = 8193.549999999999
$ws.Range("L32").Value = 15455.556
$ws.Range("M32").Value = -7906.549999999999
$ws.Range("N32").Value = -16029.556
$ws.Range("H53").Value = 23271.5
$ws.Range("I53").Value = 5000
$ws.Range("K53").Value = 5000
$ws.Range("M53").Value = -4318
$ws.Range("H61").Value = 31252406
$ws.Range("I61").Value = 41669040
$ws.Range("J61").Value = 2506
$ws.Range("K61").Value = 41669040
$ws.Range("L61").Value = 2506
$ws.Range("M61").Value = -41668828
$ws.Range("N61").Value = -2930
$ws.Range("H97").Value = 9402.5
$ws.Range("I97").Value = 12271
$ws.Range("J97").Value = 797
$ws.Range("K97").Value = 12271
$ws.Range("L97").Value = 797
$ws.Range("M97").Value = -11775
$ws.Range("N97").Value = -1789
$ws.Range("H132").Value = 8622870
$ws.Range("I132").Value = 11365620
$ws.Range("J132").Value = 2798.5715
$ws.Range("K132").Value = 34096860
$ws.Range("L132").Value = 8395.7145
$ws.Range("M132").Value = -34094330
$ws.Range("N132").Value = -13455.7145
$ws.Range("H136").Value = 31252406
$ws.Range("I136").Value = 41669040
$ws.Range("J136").Value = 2506
$ws.Range("K136").Value = 125007120
$ws.Range("L136").Value = 7518
$ws.Range("M136").Value = -125004570
$ws.Range("N136").Value = -12618

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
$ws.Range("H134").Value = 3632.1875
$ws.Range("I134").Value = 1999.4783
$ws.Range("J134").Value = 7804.6665
$ws.Range("K134").Value = 5998.4349
$ws.Range("L134").Value = 23413.9995
$ws.Range("M134").Value = -3463.4349
$ws.Range("N134").Value = -28483.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4832
$ws.Range("I58").Value = 930.6667
$ws.Range("K58").Value = 930.6667
$ws.Range("M58").Value = -727.6667
$ws.Range("H62").Value = 2850
$ws.Range("I62").Value = 2350
$ws.Range("J62").Value = 3350
$ws.Range("K62").Value = 2350
$ws.Range("L62").Value = 3350
$ws.Range("M62").Value = -1726
$ws.Range("N62").Value = -4598
$ws.Range("H65").Value = 2850
$ws.Range("I65").Value = 2350
$ws.Range("J65").Value = 3350
$ws.Range("K65").Value = 11750
$ws.Range("L65").Value = 16750
$ws.Range("M65").Value = -8630
$ws.Range("N65").Value = -22990
$ws.Range("H86").Value = 4582
$ws.Range("I86").Value = 5200
$ws.Range("J86").Value = 4170
$ws.Range("K86").Value = 5200
$ws.Range("L86").Value = 4170
$ws.Range("M86").Value = -4077
$ws.Range("N86").Value = -6416
$ws.Range("H89").Value = 4582
$ws.Range("I89").Value = 5200
$ws.Range("J89").Value = 4170
$ws.Range("K89").Value = 26000
$ws.Range("L89").Value = 20850
$ws.Range("M89").Value = -20384
$ws.Range("N89").Value = -32082
$ws.Range("H136").Value = 4832
$ws.Range("I136").Value = 930.6667
$ws.Range("K136").Value = 2792.0001
$ws.Range("M136").Value = -242.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 926924.4399999999
$ws.Range("I2").Value = 121.9
$ws.Range("K2").Value = 731.4000000000001
$ws.Range("M2").Value = -618.4000000000001
$ws.Range("H26").Value = 351.53333
$ws.Range("I26").Value = 149.3
$ws.Range("J26").Value = 756
$ws.Range("K26").Value = 447.9
$ws.Range("L26").Value = 2268
$ws.Range("M26").Value = -159.9
$ws.Range("N26").Value = -2844
$ws.Range("H118").Value = 1902.0416
$ws.Range("J118").Value = 1825.1111
$ws.Range("L118").Value = 5475.3333
$ws.Range("N118").Value = -7961.3333
$ws.Range("H131").Value = 1068.64
$ws.Range("I131").Value = 620.8333
$ws.Range("J131").Value = 1210.0526
$ws.Range("K131").Value = 1862.4999
$ws.Range("L131").Value = 3630.1578
$ws.Range("M131").Value = 3177.5001
$ws.Range("N131").Value = -13710.1578
$ws.Range("H132").Value = 1346.7778
$ws.Range("I132").Value = 684.2
$ws.Range("J132").Value = 2175
$ws.Range("K132").Value = 6157.8
$ws.Range("L132").Value = 19575
$ws.Range("M132").Value = -3627.8
$ws.Range("N132").Value = -24635
$ws.Range("H134").Value = 3349.6765
$ws.Range("I134").Value = 1603.8695
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 4811.6085
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = 258.3914999999997
$ws.Range("N134").Value = -31140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2668390.8
$ws.Range("I122").Value = 3704914.5
$ws.Range("K122").Value = 11114743.5
$ws.Range("M122").Value = -11112293.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5331.3667
$ws.Range("I7").Value = 6058.8184
$ws.Range("J7").Value = 4910.2104
$ws.Range("K7").Value = 6058.8184
$ws.Range("L7").Value = 4910.2104
$ws.Range("M7").Value = -5946.8184
$ws.Range("N7").Value = -5134.2104
$ws.Range("H22").Value = 1054.7632
$ws.Range("I22").Value = 498
$ws.Range("J22").Value = 1139.1212
$ws.Range("K22").Value = 498
$ws.Range("L22").Value = 1139.1212
$ws.Range("M22").Value = -203
$ws.Range("N22").Value = -1729.1212
$ws.Range("H27").Value = 1054.7632
$ws.Range("I27").Value = 498
$ws.Range("J27").Value = 1139.1212
$ws.Range("K27").Value = 498
$ws.Range("L27").Value = 1139.1212
$ws.Range("M27").Value = -391
$ws.Range("N27").Value = -1353.1212
$ws.Range("H68").Value = 1680.9474
$ws.Range("I68").Value = 1693.4286
$ws.Range("J68").Value = 1646
$ws.Range("K68").Value = 1693.4286
$ws.Range("L68").Value = 1646
$ws.Range("M68").Value = -944.4286
$ws.Range("N68").Value = -3144
$ws.Range("H71").Value = 1680.9474
$ws.Range("I71").Value = 1693.4286
$ws.Range("J71").Value = 1646
$ws.Range("K71").Value = 8467.143
$ws.Range("L71").Value = 8230
$ws.Range("M71").Value = -4723.143
$ws.Range("N71").Value = -15718
$ws.Range("H126").Value = 5331.3667
$ws.Range("I126").Value = 6058.8184
$ws.Range("J126").Value = 4910.2104
$ws.Range("K126").Value = 18176.4552
$ws.Range("L126").Value = 14730.6312
$ws.Range("M126").Value = -15706.4552
$ws.Range("N126").Value = -19670.6312
$ws.Range("H140").Value = 51857.43
$ws.Range("J140").Value = 51857.43
$ws.Range("L140").Value = 51857.43
$ws.Range("N140").Value = -62217.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2402.9092
$ws.Range("I126").Value = 1621.1177
$ws.Range("K126").Value = 4863.3531
$ws.Range("M126").Value = -2393.3531
